$wb = $excel.ActiveWorkbook
$flags = $wb.Worksheets.Item("Flags")
$tests = $wb.Worksheets.Item("Tests")

$flags.Range("B3").Value = "Debug"
$flags.Range("B4").Value = "False"

$tests.Range("B42").ClearContents()
$tests.Range("C42").ClearContents()
$tests.Range("D42").ClearContents()
